$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (Team1): B2=1, C2=1, D2=0, F2=1 (E2 stays empty)
$ws.Range("B2").Value = 1
$ws.Range("C2").Value = 1
$ws.Range("D2").Value = 0
$ws.Range("F2").Value = 1

# Row 4 (Team3): B4=1, D4:U4 = 1 (C4 stays empty)
$ws.Range("B4").Value = 1
$ws.Range("D4:U4").Value = 1

# Move the selection cursor to M12 (matches the saved selection in the diff)
$ws.Range("M12").Select()
